# Update "want to go" (想去人数) counts in column F across the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets,
# reflecting a refreshed data scrape.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1475
$wsExpo.Range("F7").Value = 140
$wsExpo.Range("F8").Value = 6192
$wsExpo.Range("F12").Value = 5083
$wsExpo.Range("F14").Value = 177
$wsExpo.Range("F15").Value = 1173
$wsExpo.Range("F19").Value = 6
$wsExpo.Range("F22").Value = 3595

# --- Sheet "演出" ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 75

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 75
$wsAll.Range("F5").Value = 1475
$wsAll.Range("F8").Value = 140
$wsAll.Range("F9").Value = 6192
$wsAll.Range("F13").Value = 5083
$wsAll.Range("F15").Value = 177
$wsAll.Range("F16").Value = 1173
$wsAll.Range("F20").Value = 6
$wsAll.Range("F23").Value = 3595
